# Updates the betting-odds grid (columns F:AO, rows 2:11) on Sheet1 with
# the refreshed values for Jogos_do_Dia_Betfair_Back_Lay_2026-01-01.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each array holds the new values for columns F..AO (36 columns) for the row.
$row2  = @(1.83, 1.89, 4.3, 4.6, 4.2, 4.4, 1.5, 1.04, 5.1, 1.23, 2.26, 1.75, 1.52, 2.86, 1.69, 2.3, 1.27, 2.12, 24, 24, 1000, 1000, 13, 10.5, 20, 1000, 14, 10.5, 21, 1000, 24, 20, 27, 1000, 11, 1000)
$row3  = @(2.36, 11, 1.17, 3.65, 3.15, 3.5, 1.58, 1.11, 3, 1.46, 1.66, 2.46, 1.24, 4.7, 2.1, 1.85, 1.01, 1.01, 11, 11.5, 28, 85, 8.8, 7.6, 18.5, 65, 17.5, 14.5, 22, 70, 44, 40, 55, 190, 38, 75)
$row4  = @(1.68, 1.71, 6, 6.6, 3.85, 4.1, 1.52, 1.08, 3.35, 1.4, 1.65, 2.38, 1.29, 4.1, 2.22, 1.68, 1.17, 2.38, 13, 240, 60, 230, 6.6, 990, 30, 120, 9.4, 10.5, 32, 140, 21, 25, 48, 210, 13.5, 200)
$row5  = @(9.4, 11.5, 1.33, 1.36, 6.2, 7.2, 1.26, 1.03, 6.2, 1.16, 2.74, 1.51, 1.73, 2.3, 1.91, 1.94, 3.75, 1.1, 32, 11.5, 9.8, 12, 46, 15, 10.5, 14, 110, 36, 38, 32, 420, 160, 120, 130, 150, 4.7)
$row6  = @(2.5, 2.72, 3.15, 3.4, 3.1, 3.3, 1.47, 1.1, 3.3, 1.39, 1.8, 2.22, 1.3, 4, 1.82, 1.98, 1.42, 1.58, 23, 23, 1000, 1000, 20, 14, 24, 1000, 1000, 20, 38, 1000, 900, 1000, 1000, 1000, 1000, 1000)
$row7  = @(2.22, 2.24, 3.7, 3.8, 3.5, 3.55, 1.47, 1.08, 3.65, 1.36, 1.88, 2.1, 1.33, 3.85, 1.86, 2.1, 1.35, 1.8, 13, 13.5, 24, 70, 9, 7.4, 15.5, 42, 13, 11, 18, 55, 27, 23, 38, 100, 18.5, 50)
$row8  = @(1.54, 1.56, 7, 7.2, 4.7, 4.8, 1.34, 1.04, 4.9, 1.25, 2.3, 1.75, 1.5, 2.96, 1.88, 2.1, 1.16, 2.8, 21, 25, 60, 180, 9.2, 9.8, 24, 80, 9.2, 9.6, 21, 75, 13.5, 15, 29, 90, 7.4, 120)
$row9  = @(1.94, 2, 3.85, 4.3, 3.9, 4.4, 1.31, 1.04, 4.9, 1.22, 2.42, 1.66, 1.57, 2.64, 1.64, 2.26, 1.31, 2, 22, 19.5, 32, 120, 12.5, 10, 17, 120, 14.5, 11.5, 17, 120, 23, 19.5, 30, 200, 10.5, 42)
$row10 = @(2.42, 2.44, 3.1, 3.15, 3.7, 3.75, 1.42, 1.07, 4, 1.31, 2.04, 1.94, 1.4, 3.4, 1.76, 2.28, 1.46, 1.69, 15.5, 13, 21, 50, 11, 8, 13, 34, 15.5, 11, 16.5, 42, 34, 24, 36, 80, 18.5, 29)
$row11 = @(8.4, 8.6, 1.45, 1.46, 5.1, 5.2, 1.35, 1.05, 4.9, 1.24, 2.34, 1.73, 1.52, 2.86, 1.97, 2, 3.15, 1.13, 21, 8.8, 8.6, 12, 29, 11, 9.6, 14, 75, 30, 24, 32, 260, 120, 110, 130, 160, 6.4)

$rows = @{ 2 = $row2; 3 = $row3; 4 = $row4; 5 = $row5; 6 = $row6; 7 = $row7; 8 = $row8; 9 = $row9; 10 = $row10; 11 = $row11 }

# Column F is spreadsheet column index 6.
$firstCol = 6

foreach ($r in $rows.Keys) {
    $values = $rows[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $c = $firstCol + $i
        $ws.Cells.Item($r, $c).Value = $values[$i]
    }
}
